$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings stored as plain text (t="inlineStr") in the
# source workbook, even when they look like a simple decimal number (e.g.
# "214.90"). A bare COM .Value assignment would let Excel auto-detect such
# strings as numbers and silently drop significant trailing zeros, so for
# those specific cells we briefly force a Text number format, assign the
# value, then clear the format again so the cell keeps its original default
# (unstyled) look -- only the text content changes, same as the source diff.
$ws.Range("D2").Value = '27.083.61'
$ws.Range("E2").Value = '  +3.10%  '
$ws.Range("D3").Value = '1.655.88'
$ws.Range("E3").Value = '  +3.68%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.90'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +1.63%  '
$ws.Range("E9").Value = '  +1.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.65'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.38%  '
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("D12").Value = '1.888.50'
$ws.Range("E12").Value = '  +3.77%  '
$ws.Range("D13").Value = '1.652.37'
$ws.Range("E13").Value = '  +3.14%  '
$ws.Range("E14").Value = '  +1.93%  '
$ws.Range("E15").Value = '  +3.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.91'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.85%  '
$ws.Range("D17").Value = '27.057.49'
$ws.Range("E17").Value = '  +3.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '238.17'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.89'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.81%  '
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.44'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.71%  '
$ws.Range("E23").Value = '  +4.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.27'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.83'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  +1.93%  '
$ws.Range("E28").Value = '  +1.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.85'
$ws.Range("D29").ClearFormats()
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("E31").Value = '  +1.50%  '
$ws.Range("D32").Value = '1.528.55'
$ws.Range("E32").Value = '  +4.19%  '
$ws.Range("E33").Value = '  +2.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.05'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.59'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +8.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.41'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.577'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.886'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +7.89%  '
$ws.Range("E39").Value = '  +2.58%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '66.39'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +9.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.26'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.69%  '
$ws.Range("D44").Value = '1.795.74'
$ws.Range("E44").Value = '  +3.63%  '
$ws.Range("E45").Value = '  +2.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.921'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.14'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.76%  '
$ws.Range("D48").Value = '0.0₆0106'
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("E49").Value = '  +3.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0505'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0977'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.01%  '
